$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet "API Login with Invalid Value" right after
#    "Create Product Success" (i.e. it becomes the 2nd sheet).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "API Login with Invalid Value"
$newSheet.Move($wb.Worksheets.Item("API Create Product Success"))

$ws = $wb.Worksheets.Item("API Login with Invalid Value")

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Provider"
$ws.Range("D1").Value = "Note"
$ws.Range("E1").Value = "Message"

# ---------------------------------------------------------------------------
# 3. Data rows (2-8). Columns B (Password) hold numeric-looking text such as
#    "0501"/"1234" that must stay text so the leading zero survives.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "TotoroCarbon@gmail.com"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0501"
$ws.Range("C2").Value = "native"
$ws.Range("D2").Value = "Wrong email"
$ws.Range("E2").Value = "Login Failed"

$ws.Range("A3").Value = "Totoro@gmail.com"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1234"
$ws.Range("C3").Value = "native"
$ws.Range("D3").Value = "Wrong password"
$ws.Range("E3").Value = "Login Failed"

$ws.Range("A4").Value = "Totoro@gmail.com"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "0501"
$ws.Range("C4").Value = "facebook"
$ws.Range("D4").Value = "Wrong Provider"
$ws.Range("E4").Value = "Request Error: access token is required."

$ws.Range("C5").Value = "native"
$ws.Range("D5").Value = "Empty case 1"
$ws.Range("E5").Value = "Email and password are required."

$ws.Range("A6").Value = "Totoro@gmail.com"
$ws.Range("C6").Value = "native"
$ws.Range("D6").Value = "Empty case 2"
$ws.Range("E6").Value = "Email and password are required."

$ws.Range("A7").Value = "Totoro@gmail.com"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "0501"
$ws.Range("D7").Value = "Empty case 3"
$ws.Range("E7").Value = 'Provider can be "native" or "facebook" only'

$ws.Range("D8").Value = "Empty case 4"
$ws.Range("E8").Value = 'Provider can be "native" or "facebook" only'

# ---------------------------------------------------------------------------
# 4. Hyperlink the e-mail addresses in column A (mailto: links)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:TotoroCarbon@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:Totoro@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:Totoro@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:Totoro@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:Totoro@gmail.com")

# ---------------------------------------------------------------------------
# 5. Column widths / page setup to resemble the authored sheet
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.42
$ws.Columns.Item(2).ColumnWidth = 9.92
$ws.Columns.Item(3).ColumnWidth = 9.92
$ws.Columns.Item(4).ColumnWidth = 15.58
$ws.Columns.Item(5).ColumnWidth = 38.75
$ws.Columns.Item(6).ColumnWidth = 35.25

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 6. Make the new sheet active / selected (matches workbookView activeTab,
#    the removal of tabSelected from the first sheet, and the zoom level).
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 239
$ws.Range("E8").Select()
